$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Cells.Item(6, 8).Value = 236.44444
$ws.Cells.Item(6, 9).Value = 162.11765
$ws.Cells.Item(6, 11).Value = 486.35295
$ws.Cells.Item(6, 13).Value = -374.35295

$ws.Cells.Item(11, 8).Value = 2606.1667
$ws.Cells.Item(11, 9).Value = 2606.1667
$ws.Cells.Item(11, 11).Value = 2606.1667
$ws.Cells.Item(11, 13).Value = -2466.1667

$ws.Cells.Item(12, 8).Value = 323.33334
$ws.Cells.Item(12, 9).Value = 288.2
$ws.Cells.Item(12, 10).Value = 499
$ws.Cells.Item(12, 11).Value = 288.2
$ws.Cells.Item(12, 12).Value = 499
$ws.Cells.Item(12, 13).Value = -118.2
$ws.Cells.Item(12, 14).Value = -839

$ws.Cells.Item(19, 8).Value = 1337.2222
$ws.Cells.Item(19, 9).Value = 1318
$ws.Cells.Item(19, 11).Value = 1318
$ws.Cells.Item(19, 13).Value = -1143

$ws.Cells.Item(43, 8).Value = 3911.7646
$ws.Cells.Item(43, 9).Value = 2640.2
$ws.Cells.Item(43, 11).Value = 2640.2
$ws.Cells.Item(43, 13).Value = -2571.2

$ws.Cells.Item(53, 8).Value = 548.875
$ws.Cells.Item(53, 10).Value = 1264.6666
$ws.Cells.Item(53, 12).Value = 1264.6666
$ws.Cells.Item(53, 14).Value = -2538.6666

$ws.Cells.Item(58, 8).Value = 1175.625
$ws.Cells.Item(58, 9).Value = 1200.7142
$ws.Cells.Item(58, 10).Value = 1000
$ws.Cells.Item(58, 11).Value = 3602.1426
$ws.Cells.Item(58, 12).Value = 3000
$ws.Cells.Item(58, 13).Value = -3452.1426
$ws.Cells.Item(58, 14).Value = -3300

$ws.Cells.Item(99, 8).Value = 2743.6667
$ws.Cells.Item(99, 10).Value = 4024.5
$ws.Cells.Item(99, 12).Value = 12073.5
$ws.Cells.Item(99, 14).Value = -15069.5

$ws.Cells.Item(132, 8).Value = 7135.25
$ws.Cells.Item(132, 9).Value = 7877.4
$ws.Cells.Item(132, 11).Value = 23632.2
$ws.Cells.Item(132, 13).Value = -21102.2

$ws.Cells.Item(138, 8).Value = 2730.3235
$ws.Cells.Item(138, 9).Value = 2453.4
$ws.Cells.Item(138, 11).Value = 7360.200000000001
$ws.Cells.Item(138, 13).Value = -2220.200000000001

$ws.Cells.Item(141, 8).Value = 2179.842
$ws.Cells.Item(141, 9).Value = 1977.5883
$ws.Cells.Item(141, 11).Value = 5932.7649
$ws.Cells.Item(141, 13).Value = -752.7649000000001

$ws = $wb.Worksheets.Item("ARM")
$ws.Cells.Item(32, 8).Value = 29845.81
$ws.Cells.Item(32, 9).Value = 30641.03
$ws.Cells.Item(32, 11).Value = 30641.03
$ws.Cells.Item(32, 13).Value = -30354.03

$ws.Cells.Item(61, 8).Value = 1967.25
$ws.Cells.Item(61, 9).Value = 1967.25
$ws.Cells.Item(61, 11).Value = 1967.25
$ws.Cells.Item(61, 13).Value = -1755.25

$ws.Cells.Item(74, 8).Value = 3656.5652
$ws.Cells.Item(74, 9).Value = 3657.2273
$ws.Cells.Item(74, 11).Value = 3657.2273
$ws.Cells.Item(74, 13).Value = -2783.2273

$ws.Cells.Item(77, 8).Value = 3656.5652
$ws.Cells.Item(77, 9).Value = 3657.2273
$ws.Cells.Item(77, 11).Value = 18286.1365
$ws.Cells.Item(77, 13).Value = -13918.1365

$ws.Cells.Item(102, 8).Value = 2494.2
$ws.Cells.Item(102, 9).Value = 2007.8667
$ws.Cells.Item(102, 10).Value = 3953.2
$ws.Cells.Item(102, 11).Value = 2007.8667
$ws.Cells.Item(102, 12).Value = 3953.2
$ws.Cells.Item(102, 13).Value = -385.8667
$ws.Cells.Item(102, 14).Value = -7197.2

$ws.Cells.Item(122, 8).Value = 2379.1667
$ws.Cells.Item(122, 9).Value = 2044.3334
$ws.Cells.Item(122, 10).Value = 3048.8333
$ws.Cells.Item(122, 11).Value = 6133.0002
$ws.Cells.Item(122, 12).Value = 9146.499899999999
$ws.Cells.Item(122, 13).Value = -3683.0002
$ws.Cells.Item(122, 14).Value = -14046.4999

$ws.Cells.Item(132, 8).Value = 129124.875
$ws.Cells.Item(132, 9).Value = 171249.83
$ws.Cells.Item(132, 11).Value = 513749.49
$ws.Cells.Item(132, 13).Value = -511219.49

$ws.Cells.Item(136, 8).Value = 1967.25
$ws.Cells.Item(136, 9).Value = 1967.25
$ws.Cells.Item(136, 11).Value = 5901.75
$ws.Cells.Item(136, 13).Value = -3351.75

$ws = $wb.Worksheets.Item("BSM")
$ws.Cells.Item(96, 8).Value = 24667
$ws.Cells.Item(96, 9).Value = 19500.5
$ws.Cells.Item(96, 10).Value = 35000
$ws.Cells.Item(96, 11).Value = 19500.5
$ws.Cells.Item(96, 12).Value = 35000
$ws.Cells.Item(96, 13).Value = -16754.5
$ws.Cells.Item(96, 14).Value = -40492

$ws.Cells.Item(134, 8).Value = 2347.7812
$ws.Cells.Item(134, 9).Value = 2197.7097
$ws.Cells.Item(134, 11).Value = 6593.1291
$ws.Cells.Item(134, 13).Value = -4058.1291

$ws = $wb.Worksheets.Item("CRP")
$ws.Cells.Item(58, 8).Value = 69452.47
$ws.Cells.Item(58, 9).Value = 93163.45
$ws.Cells.Item(58, 11).Value = 93163.45
$ws.Cells.Item(58, 13).Value = -92960.45

$ws.Cells.Item(62, 8).Value = 4327.857

$ws.Cells.Item(65, 8).Value = 4327.857

$ws.Cells.Item(132, 8).Value = 1950.5
$ws.Cells.Item(132, 9).Value = 1950.5
$ws.Cells.Item(132, 11).Value = 5851.5
$ws.Cells.Item(132, 13).Value = -3321.5

$ws.Cells.Item(136, 8).Value = 69452.47
$ws.Cells.Item(136, 9).Value = 93163.45
$ws.Cells.Item(136, 11).Value = 279490.35
$ws.Cells.Item(136, 13).Value = -276940.35

$ws = $wb.Worksheets.Item("CUL")
$ws.Cells.Item(4, 8).Value = 414357
$ws.Cells.Item(4, 9).Value = 414357
$ws.Cells.Item(4, 11).Value = 1243071
$ws.Cells.Item(4, 13).Value = -1242959

$ws.Cells.Item(8, 8).Value = 495.53333
$ws.Cells.Item(8, 9).Value = 495.53333
$ws.Cells.Item(8, 11).Value = 1486.59999
$ws.Cells.Item(8, 13).Value = -1347.59999

$ws.Cells.Item(12, 8).Value = 108.27273
$ws.Cells.Item(12, 9).Value = 250
$ws.Cells.Item(12, 10).Value = 94.09999999999999
$ws.Cells.Item(12, 11).Value = 750
$ws.Cells.Item(12, 12).Value = 282.3
$ws.Cells.Item(12, 13).Value = -577
$ws.Cells.Item(12, 14).Value = -628.3

$ws.Cells.Item(19, 8).Value = 7305.6
$ws.Cells.Item(19, 9).Value = 28
$ws.Cells.Item(19, 10).Value = 9125
$ws.Cells.Item(19, 11).Value = 84
$ws.Cells.Item(19, 12).Value = 27375
$ws.Cells.Item(19, 13).Value = 90
$ws.Cells.Item(19, 14).Value = -27723

$ws.Cells.Item(35, 8).Value = 13333.333
$ws.Cells.Item(35, 9).Value = 15000
$ws.Cells.Item(35, 10).Value = 12500
$ws.Cells.Item(35, 11).Value = 45000
$ws.Cells.Item(35, 12).Value = 37500
$ws.Cells.Item(35, 13).Value = -44712
$ws.Cells.Item(35, 14).Value = -38076

$ws.Cells.Item(47, 8).Value = 4000
$ws.Cells.Item(47, 9).Value = 4000
$ws.Cells.Item(47, 11).Value = 12000
$ws.Cells.Item(47, 13).Value = -11569

$ws.Cells.Item(57, 8).Value = 13230
$ws.Cells.Item(57, 9).Value = 6460
$ws.Cells.Item(57, 11).Value = 19380
$ws.Cells.Item(57, 13).Value = -18821

$ws.Cells.Item(64, 8).Value = 6999
$ws.Cells.Item(64, 9).Value = 6999
$ws.Cells.Item(64, 11).Value = 20997
$ws.Cells.Item(64, 13).Value = -20727

$ws.Cells.Item(67, 8).Value = 6999
$ws.Cells.Item(67, 9).Value = 6999
$ws.Cells.Item(67, 11).Value = 20997
$ws.Cells.Item(67, 13).Value = -20061

$ws.Cells.Item(129, 8).Value = 2573.9285
$ws.Cells.Item(129, 9).Value = 2116.7144
$ws.Cells.Item(129, 10).Value = 3031.1428
$ws.Cells.Item(129, 11).Value = 6350.1432
$ws.Cells.Item(129, 12).Value = 9093.428400000001
$ws.Cells.Item(129, 13).Value = -1350.1432
$ws.Cells.Item(129, 14).Value = -19093.4284

$ws.Cells.Item(130, 8).Value = 250005000
$ws.Cells.Item(130, 9).Value = 500000000
$ws.Cells.Item(130, 10).Value = 10000
$ws.Cells.Item(130, 11).Value = 1500000000
$ws.Cells.Item(130, 12).Value = 30000
$ws.Cells.Item(130, 13).Value = -1499994980
$ws.Cells.Item(130, 14).Value = -40040

$ws.Cells.Item(133, 8).Value = 9064.375
$ws.Cells.Item(133, 10).Value = 25000
$ws.Cells.Item(133, 12).Value = 75000
$ws.Cells.Item(133, 14).Value = -85120

$ws.Cells.Item(137, 8).Value = 3235.1538
$ws.Cells.Item(137, 9).Value = 2158.75
$ws.Cells.Item(137, 10).Value = 4957.4
$ws.Cells.Item(137, 11).Value = 6476.25
$ws.Cells.Item(137, 12).Value = 14872.2
$ws.Cells.Item(137, 13).Value = -1376.25
$ws.Cells.Item(137, 14).Value = -25072.2

$ws = $wb.Worksheets.Item("GSM")
$ws.Cells.Item(70, 8).Value = 2213.1177
$ws.Cells.Item(70, 9).Value = 1641
$ws.Cells.Item(70, 11).Value = 1641
$ws.Cells.Item(70, 13).Value = -1371

$ws.Cells.Item(73, 8).Value = 2213.1177
$ws.Cells.Item(73, 9).Value = 1641
$ws.Cells.Item(73, 11).Value = 1641
$ws.Cells.Item(73, 13).Value = -705

$ws.Cells.Item(97, 8).Value = 1030.9584
$ws.Cells.Item(97, 9).Value = 809.13336
$ws.Cells.Item(97, 10).Value = 1400.6666
$ws.Cells.Item(97, 11).Value = 809.13336
$ws.Cells.Item(97, 12).Value = 1400.6666
$ws.Cells.Item(97, 13).Value = -313.13336
$ws.Cells.Item(97, 14).Value = -2392.6666

$ws.Cells.Item(102, 8).Value = 4534.8335
$ws.Cells.Item(102, 10).Value = 4999
$ws.Cells.Item(102, 12).Value = 4999
$ws.Cells.Item(102, 14).Value = -8243

$ws.Cells.Item(140, 8).Value = 0
$ws.Cells.Item(140, 10).Value = 0
$ws.Cells.Item(140, 14).Value = 0
$ws.Cells.Item(140, 12).ClearContents()

$ws = $wb.Worksheets.Item("LTW")
$ws.Cells.Item(25, 8).Value = 20000
$ws.Cells.Item(25, 10).Value = 20000
$ws.Cells.Item(25, 12).Value = 20000
$ws.Cells.Item(25, 14).Value = -20460

$ws.Cells.Item(40, 8).Value = 5499.231
$ws.Cells.Item(40, 9).Value = 4748.1
$ws.Cells.Item(40, 11).Value = 4748.1
$ws.Cells.Item(40, 13).Value = -4612.1

$ws.Cells.Item(93, 8).Value = 2379.4
$ws.Cells.Item(93, 9).Value = 2476.875
$ws.Cells.Item(93, 11).Value = 2476.875
$ws.Cells.Item(93, 13).Value = -1228.875

$ws.Cells.Item(100, 8).Value = 2647.9333
$ws.Cells.Item(100, 9).Value = 2714.3333
$ws.Cells.Item(100, 11).Value = 2714.3333
$ws.Cells.Item(100, 13).Value = -2173.3333

$ws.Cells.Item(136, 8).Value = 5678.4443
$ws.Cells.Item(136, 9).Value = 5228.7144
$ws.Cells.Item(136, 10).Value = 7252.5
$ws.Cells.Item(136, 11).Value = 15686.1432
$ws.Cells.Item(136, 12).Value = 21757.5
$ws.Cells.Item(136, 13).Value = -13136.1432
$ws.Cells.Item(136, 14).Value = -26857.5

$ws = $wb.Worksheets.Item("WVR")
$ws.Cells.Item(7, 8).Value = 1000
$ws.Cells.Item(7, 9).Value = 1000
$ws.Cells.Item(7, 11).Value = 1000
$ws.Cells.Item(7, 13).Value = -887

$ws.Cells.Item(13, 8).Value = 599
$ws.Cells.Item(13, 9).Value = 599
$ws.Cells.Item(13, 11).Value = 599
$ws.Cells.Item(13, 13).Value = -459

$ws.Cells.Item(81, 8).Value = 2638.5
$ws.Cells.Item(81, 9).Value = 1028.6666
$ws.Cells.Item(81, 10).Value = 4248.3335
$ws.Cells.Item(81, 11).Value = 2057.3332
$ws.Cells.Item(81, 12).Value = 8496.666999999999
$ws.Cells.Item(81, 13).Value = -996.3332
$ws.Cells.Item(81, 14).Value = -10618.667

$ws.Cells.Item(84, 8).Value = 2638.5
$ws.Cells.Item(84, 9).Value = 1028.6666
$ws.Cells.Item(84, 10).Value = 4248.3335
$ws.Cells.Item(84, 11).Value = 10286.666
$ws.Cells.Item(84, 12).Value = 42483.335
$ws.Cells.Item(84, 13).Value = -4982.666000000001
$ws.Cells.Item(84, 14).Value = -53091.335
